$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Test Case description (B1)
$ws.Range("B1").Value = "Test Case: Testing to see if creating assessments reaches data base with no errors"

# Rewrite Step 1 / Expected Result (row 2)
$ws.Range("C2").Value = "Step 1: While logged out go to the create assessment page"
$ws.Range("D2").Value = "I am returned to the login page"

# Rewrite Step 2 / Expected Result (row 3)
$ws.Range("C3").Value = "Step 2: Login as a the appropriate and go to the `"create assessment`" page"
$ws.Range("D3").Value = "I am redirected to the create assessment page"

# New Step 3 / Expected Result (row 4)
$ws.Range("C4").Value = "Step 3: Fill out a assessment for an employee"
$ws.Range("D4").Value = "A new assessment will be added to the database on that employee"

# New Step 5 / Expected Result (row 5)
$ws.Range("C5").Value = "Step 5: Leave Certain fields that are required blank (everything but comments are required)"
$ws.Range("D5").Value = "Error text pop up saying that some of the fields have been left blank"

# New Step 6 / Expected Result (row 6)
$ws.Range("C6").Value = "Step 6: Create an assessment for the current role logged in"
$ws.Range("D6").Value = "I am denied access to this"

# Clear the explicit cell selection that was saved with the original workbook
$ws.Range("A1").Select()
